# Add newly-documented classes/functions rows to the idsw workbook.
#
# 1) "idsw.datafetch.core" sheet: a new "SharePointDownloader" class with
#    five methods (get_token, get_response_id, get_drive_id, find_file,
#    download_file) appended after the existing "IngestExcelTables" rows.
# 2) "idsw.modelling.nonsupervised" sheet: a new "benford_outliers_detection"
#    function appended after the existing rows.

$wb = $excel.ActiveWorkbook

# --- idsw.datafetch.core ---------------------------------------------------
$wsCore = $wb.Worksheets.Item("idsw.datafetch.core")

$lastRow = $wsCore.Cells.Item($wsCore.Rows.Count, 1).End(-4162).Row
# lastRow currently 53 (A53 = 52). New entries start at row 54.

$newRows = @(
    @("SharePointDownloader", "get_token"),
    @("SharePointDownloader", "get_response_id"),
    @("SharePointDownloader", "get_drive_id"),
    @("SharePointDownloader", "find_file"),
    @("SharePointDownloader", "download_file")
)

$idValue = [int]$wsCore.Cells.Item($lastRow, 1).Value()
$r = $lastRow + 1
foreach ($row in $newRows) {
    $idValue = $idValue + 1
    $wsCore.Cells.Item($r, 1).Value = $idValue
    $wsCore.Cells.Item($r, 2).Value = $row[0]
    $wsCore.Cells.Item($r, 3).Value = $row[1]
    $r = $r + 1
}

# --- idsw.modelling.nonsupervised ------------------------------------------
$wsNonsup = $wb.Worksheets.Item("idsw.modelling.nonsupervised")

$lastRowNonsup = $wsNonsup.Cells.Item($wsNonsup.Rows.Count, 1).End(-4162).Row
# lastRowNonsup currently 5 (A5 = 4). New entry goes at row 6.

$rNonsup = $lastRowNonsup + 1
$idValueNonsup = [int]$wsNonsup.Cells.Item($lastRowNonsup, 1).Value() + 1
$wsNonsup.Cells.Item($rNonsup, 1).Value = $idValueNonsup
$wsNonsup.Cells.Item($rNonsup, 3).Value = "benford_outliers_detection"
